$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sample IDs in column A:
#  rows 6-19  -> S1_2017 .. S14_2017
#  rows 20-38 -> S15_2018 .. S33_2018
for ($i = 1; $i -le 14; $i++) {
    $row = 5 + $i
    $ws.Cells.Item($row, 1).Value = "S" + $i + "_2017"
}
for ($i = 15; $i -le 33; $i++) {
    $row = 5 + $i
    $ws.Cells.Item($row, 1).Value = "S" + $i + "_2018"
}

# Update the frozen-pane view / selection to match the new scroll position.
$ws.Range("C10").Select()
